$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9410267472267151
$ws.Range("B1").Value = 3.140889644622803
$ws.Range("C1").Value = 6.86278772354126
$ws.Range("D1").Value = 1.956629157066345
$ws.Range("E1").Value = 1.375646591186523
